# RollerSystem.pptx text edits
#  - Slide 2 ("Introduction"): reword the intro sentence about the roller
#    shutters and add a closing period.
#  - Slide 5 ("Functioning and Possible Extension"): reword the "lower than
#    a minimum threshold" sentence.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - Content Placeholder 2
# ---------------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$shp = $s2.Shapes.Item("Content Placeholder 2")
$tr  = $shp.TextFrame.TextRange

# 1) " on the roller " -> " on the house roller "
$full = $tr.Text
$old  = " on the roller "
$new  = " on the house roller "
$idx  = $full.IndexOf($old)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $old.Length)
    $rng.Text = $new
}

# 2) " of an house in " -> " in "
$full = $tr.Text
$old  = " of an house in "
$new  = " in "
$idx  = $full.IndexOf($old)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $old.Length)
    $rng.Text = $new
}

# 3) End the "...outside brightness" paragraph with a period.
$para1 = $tr.Paragraphs(1, 1)
$null = $para1.InsertAfter(".")

# ---------------------------------------------------------------------
# Slide 5 - Content Placeholder 2
# ---------------------------------------------------------------------
$s5  = $p.Slides.Item(5)
$shp = $s5.Shapes.Item("Content Placeholder 2")
$tr  = $shp.TextFrame.TextRange

# "Otherwise, if it is lower..." -> "Otherwise, if the value is lower..."
$para3 = $tr.Paragraphs(3, 1)
$start = $para3.Start
$len   = $para3.Length
$newText = "Otherwise, if the value is lower than a minimum threshold the roller is closed."

$null = $para3.InsertBefore($newText)
$oldRng = $tr.Characters($start + $newText.Length, $len)
$oldRng.Delete()
